$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 15) following the same pattern as row 14
$ws.Range("B15").Value = 45294
$ws.Range("B15").NumberFormat = $ws.Range("B14").NumberFormat
$ws.Range("C15").Value = 11
$ws.Range("D15").Value = 22
$ws.Range("E15").Value = 33

# Update selection to match the diff
$ws.Range("B16").Select()
